# Daily Scrum Report - add "Meeting 2" actual content and a new "Meeting 3"
# template block (interview consent form related notes were folded into
# this deliverable update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Fill in the real status-update text for "Meeting 2" (rows 35-57),
# replacing the generic placeholder headers that were still in the
# template. New shared strings are appended in sheet (row) order so the
# resulting sharedStrings.xml table matches what Excel itself would
# produce.
# ---------------------------------------------------------------------

# Arpit
$ws.Range("C35").Value = "What they have done: Absent from meeting"
$ws.Range("C37").Value = "What impediments prevent them from progressing: None"

# Brody
$ws.Range("C39").Value = "What they have done: Installed Android Studio, Reqs and Use Cases for stories 4-6"
$ws.Range("C40").Value = "What they will do: interview scripts"

# Michael
$ws.Range("C43").Value = "What they have done: Absent from meeting"
$ws.Range("C45").Value = "What impediments prevent them from progressing: "

# Sakshyam
$ws.Range("C47").Value = "What they have done: Installed Android Studio, user reqs for 1-3"
$ws.Range("C48").Value = "What they will do: system reqs and use cases for 1-3, mockup"
$ws.Range("C49").Value = "What impediments prevent them from progressing: None"

# Vasilis
$ws.Range("C51").Value = "What they have done: reqs and use cases for stories 13-15, created various documents"
$ws.Range("C52").Value = "What they will do: Permission form"
$ws.Range("C53").Value = "What impediments prevent them from progressing: None"

# Yong
$ws.Range("C55").Value = "What they have done: Installed Android Studio, reqs for 10-12"
$ws.Range("C56").Value = "What they will do: reqs for 19, use cases for 10-19 and 19"
$ws.Range("C57").Value = "What impediments prevent them from progressing: None"

# ---------------------------------------------------------------------
# Append a brand new "Meeting 3" block (rows 58-82), mirroring the same
# blank template layout used for the other meetings: a highlighted
# separator row, then one row per team member with the three standard
# prompts underneath.
# ---------------------------------------------------------------------

# Highlighted separator row, copying the look of the existing one (row 33)
$ws.Range("A33:L33").Copy($ws.Range("A58:L58"))

# Meeting 3 heading
$ws.Range("A59").Value = "Meeting 3"
$ws.Range("B59").Value = "xx Feb 18"

# Arpit
$ws.Range("B60").Value = "Arpit"
$ws.Range("C60").Value = "What they have done:"
$ws.Range("C61").Value = "What they will do:"
$ws.Range("C62").Value = "What impediments prevent them from progressing:"

# Brody
$ws.Range("B64").Value = "Brody"
$ws.Range("C64").Value = "What they have done:"
$ws.Range("C65").Value = "What they will do:"
$ws.Range("C66").Value = "What impediments prevent them from progressing:"

# Michael
$ws.Range("B68").Value = "Michael"
$ws.Range("C68").Value = "What they have done:"
$ws.Range("C69").Value = "What they will do:"
$ws.Range("C70").Value = "What impediments prevent them from progressing:"

# Sakshyam
$ws.Range("B72").Value = "Sakshyam"
$ws.Range("C72").Value = "What they have done:"
$ws.Range("C73").Value = "What they will do:"
$ws.Range("C74").Value = "What impediments prevent them from progressing:"

# Vasilis
$ws.Range("B76").Value = "Vasilis"
$ws.Range("C76").Value = "What they have done:"
$ws.Range("C77").Value = "What they will do:"
$ws.Range("C78").Value = "What impediments prevent them from progressing:"

# Yong
$ws.Range("B80").Value = "Yong"
$ws.Range("C80").Value = "What they have done:"
$ws.Range("C81").Value = "What they will do:"
$ws.Range("C82").Value = "What impediments prevent them from progressing:"

# ---------------------------------------------------------------------
# Update the saved view state to match where the editor ended up.
# ---------------------------------------------------------------------
$app = $wb.Application
$win = $app.ActiveWindow
$win.ScrollRow = 29
$win.ScrollColumn = 1
$ws.Range("H47").Select()
